# Generate Report for Handback
# Updates the localization-status workbook: marks the two files as handed
# back (status + target/handback columns) and widens a few columns that
# now hold longer content.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$c0c2Name = "c0c2e956-77ee-48db-9b90-5d8caf543670.md"
$da29Name = "da291fba-1a87-4f24-8484-531a5aaaa0b1.md"

$c0c2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b571e0ad554d6780db8549f66c4ce7a29726199c/e2e/c0c2e956-77ee-48db-9b90-5d8caf543670.md"
$da29Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b571e0ad554d6780db8549f66c4ce7a29726199c/e2e/da291fba-1a87-4f24-8484-531a5aaaa0b1.md"

# Column widths (character units). This runtime rounds ColumnWidth to the
# nearest pixel internally, so these are the closest achievable settings
# to the new, much wider columns used for the long file names below.
$wideWidth = 29.166666666666668
$fullWidth = 39.166666666666664

# ---------------------------------------------------------------------
# Overview sheet: both language status columns flip to "Handed back"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Range("E3").Value = $statusText
$ws1.Range("F3").Value = $statusText

$ws1.Columns.Item(5).ColumnWidth = $wideWidth
$ws1.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# Helper routine applied to both the zh-cn and de-de language sheets.
# $handback2 / $handback3 are the language-specific "latest handback
# file" names; $datetimeText is the language-specific handback timestamp.
# ---------------------------------------------------------------------
function Update-LanguageSheet($ws, $handback2, $handback3, $datetimeText) {
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    $ws.Range("I2").Value = $c0c2Name
    $ws.Range("J2").Value = $handback2
    $ws.Range("K2").Value = $datetimeText

    $ws.Range("I3").Value = $da29Name
    $ws.Range("J3").Value = $handback3
    $ws.Range("K3").Value = $datetimeText

    # Re-create the hyperlinks so the relationship ids line up as
    # A2, I2, A3, I3 (matching the order new links were appended).
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $c0c2Url, "", "", $c0c2Name)
    $ws.Hyperlinks.Add($ws.Range("I2"), $c0c2Url, "", "", $c0c2Name)
    $ws.Hyperlinks.Add($ws.Range("A3"), $da29Url, "", "", $da29Name)
    $ws.Hyperlinks.Add($ws.Range("I3"), $da29Url, "", "", $da29Name)

    $ws.Columns.Item(3).ColumnWidth = $wideWidth
    $ws.Columns.Item(9).ColumnWidth = $fullWidth
    $ws.Columns.Item(10).ColumnWidth = $fullWidth
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$zhHandback2 = "c0c2e956-77ee-48db-9b90-5d8caf543670.f1eb9b27dfabe7e8315b4ffb98a70850414c201b.zh-cn.xlf"
$zhHandback3 = "da291fba-1a87-4f24-8484-531a5aaaa0b1.a7a1849f51c6c31be3a893e45b58671191a2a4c9.zh-cn.xlf"
Update-LanguageSheet $ws2 $zhHandback2 $zhHandback3 "2016-10-13 14:04:36"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$deHandback2 = "c0c2e956-77ee-48db-9b90-5d8caf543670.f1eb9b27dfabe7e8315b4ffb98a70850414c201b.de-de.xlf"
$deHandback3 = "da291fba-1a87-4f24-8484-531a5aaaa0b1.a7a1849f51c6c31be3a893e45b58671191a2a4c9.de-de.xlf"
Update-LanguageSheet $ws3 $deHandback2 $deHandback3 "2016-10-13 14:04:52"

Write-Host "Handback report generated."
